$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended after existing row 33 (regcntr_id 10005 block)
$newRows = @(
    @(10005, 110033, 10005),
    @(10005, 110034, 10005),
    @(10005, 110035, 10005)
)

$r = 34
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
    $r++
}

$ws.Range($ws.Rows.Item(37), $ws.Rows.Item($ws.Rows.Count)).Select()

$wb.Save()
